# Insert a new data row at row 595 (pushing the existing rows 595:670 down
# to 596:671) and populate it with the new weekly price record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("595").Insert()

$ws.Cells.Item(595, 1).Value = 8
$ws.Cells.Item(595, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(595, 3).Value = "Coquimbo"
$ws.Cells.Item(595, 4).Value = 45142
$ws.Cells.Item(595, 5).Value = 4
$ws.Cells.Item(595, 6).Value = 100112043
$ws.Cells.Item(595, 7).Value = "Pepino dulce"
$ws.Cells.Item(595, 8).Value = "Sin especificar"
$ws.Cells.Item(595, 9).Value = "Primera"
$ws.Cells.Item(595, 10).Value = 200
$ws.Cells.Item(595, 11).Value = 16000
$ws.Cells.Item(595, 12).Value = 17000
$ws.Cells.Item(595, 13).Value = 16500
$ws.Cells.Item(595, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(595, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(595, 16).Value = 917
$ws.Cells.Item(595, 17).Value = 18
$ws.Cells.Item(595, 18).Value = "Hortaliza"
